$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.196.35"
$ws.Range("E2").Value = "'  +0.38%  "
$ws.Range("D3").Value = "'1.904.26"
$ws.Range("E3").Value = "'  +0.86%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'306.54"
$ws.Range("E5").Value = "'  -0.29%  "
$ws.Range("E6").Value = "'  +0.02%  "
$ws.Range("D7").Value = "'0.5264"
$ws.Range("E7").Value = "'  +2.19%  "
$ws.Range("E8").Value = "'  +1.61%  "
$ws.Range("D9").Value = "'0.07253"
$ws.Range("E9").Value = "'  +0.56%  "
$ws.Range("E10").Value = "'  +0.62%  "
$ws.Range("D11").Value = "'0.9003"
$ws.Range("E11").Value = "'  -0.37%  "
$ws.Range("D12").Value = "'0.08380"
$ws.Range("E12").Value = "'  +10.00%  "
$ws.Range("D13").Value = "'1.906.86"
$ws.Range("E13").Value = "'  +0.91%  "
$ws.Range("D14").Value = "'94.95"
$ws.Range("E14").Value = "'  +0.23%  "
$ws.Range("D15").Value = "'5.280"
$ws.Range("E15").Value = "'  +0.12%  "
$ws.Range("D16").Value = "'0.9977"
$ws.Range("E16").Value = "'  -0.26%  "
$ws.Range("D17").Value = "'0.000008609"
$ws.Range("D18").Value = "'14.56"
$ws.Range("E18").Value = "'  +1.40%  "
$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "'  +0.00%  "
$ws.Range("D20").Value = "'27.234.52"
$ws.Range("E20").Value = "'  +0.36%  "
$ws.Range("E21").Value = "'  +0.27%  "
$ws.Range("D22").Value = "'2.151.66"
$ws.Range("E22").Value = "'  +0.57%  "
$ws.Range("E23").Value = "'  +0.53%  "
$ws.Range("D24").Value = "'6.443"
$ws.Range("E24").Value = "'  +0.27%  "
$ws.Range("D25").Value = "'147.43"
$ws.Range("E25").Value = "'  +1.20%  "
$ws.Range("D26").Value = "'2.284"
$ws.Range("E26").Value = "'  +5.03%  "
$ws.Range("D27").Value = "'1.753"
$ws.Range("E27").Value = "'  -2.20%  "
$ws.Range("D28").Value = "'18.18"
$ws.Range("E28").Value = "'  +0.73%  "
$ws.Range("D29").Value = "'114.87"
$ws.Range("E29").Value = "'  +0.23%  "
$ws.Range("D30").Value = "'4.928"
$ws.Range("E30").Value = "'  -1.01%  "
$ws.Range("D31").Value = "'4.822"
$ws.Range("E31").Value = "'  -0.02%  "
$ws.Range("D32").Value = "'0.09290"
$ws.Range("E32").Value = "'  +0.90%  "
$ws.Range("D33").Value = "'0.8101"
$ws.Range("E33").Value = "'  +6.76%  "
$ws.Range("E34").Value = "'  -0.06%  "
$ws.Range("D35").Value = "'1.241"
$ws.Range("E35").Value = "'  +3.53%  "
$ws.Range("D36").Value = "'2.952"
$ws.Range("E36").Value = "'  -1.95%  "
$ws.Range("D37").Value = "'3.383"
$ws.Range("E37").Value = "'  +3.41%  "
$ws.Range("D38").Value = "'2.624"
$ws.Range("E38").Value = "'  +2.34%  "
$ws.Range("D39").Value = "'0.5727"
$ws.Range("E39").Value = "'  +1.82%  "
$ws.Range("D40").Value = "'0.01987"
$ws.Range("E40").Value = "'  -0.27%  "
$ws.Range("E41").Value = "'  -0.10%  "
$ws.Range("D42").Value = "'6.646"
$ws.Range("E42").Value = "'  +1.09%  "
$ws.Range("D43").Value = "'8.979"
$ws.Range("E43").Value = "'  -0.60%  "
$ws.Range("D44").Value = "'117.52"
$ws.Range("E44").Value = "'  -0.71%  "
$ws.Range("D45").Value = "'0.1513"
$ws.Range("E45").Value = "'  +0.42%  "
$ws.Range("D46").Value = "'0.4850"
$ws.Range("E46").Value = "'  +1.04%  "
$ws.Range("B47").Value = "'PaxDollar"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.0000"
$ws.Range("E47").Value = "'  +0.06%  "
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.12"
$ws.Range("E48").Value = "'  -0.69%  "
$ws.Range("D49").Value = "'1.616"
$ws.Range("E49").Value = "'  +2.50%  "
$ws.Range("D50").Value = "'37.48"
$ws.Range("E50").Value = "'  +0.79%  "
$ws.Range("D51").Value = "'63.79"
$ws.Range("E51").Value = "'  +0.30%  "
